# Saldo_guide.xlsx update
# - Rename sheet to reflect the new extraction timestamp (20240913-104913 -> 20240916-100029)
# - Bump "Dt. Referencia" (column G) for every data row from 45548 (2024-09-13) to 45551 (2024-09-16)
# - Refresh "Vl. Projetado" (D), "Saldo Previsto" (E) and "Vl. Total" (H) for the accounts whose
#   projected/pending values were settled (or replaced by a new pending value) since the last extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet (tab) name -> new export timestamp
$ws.Name = "IClientBalance-20240916-100029-"

# 2) New reference date for every data row (row 2 .. 274), column G
$newDate = 45551
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = $newDate
}

# 3) Rows whose Vl. Projetado (D) / Saldo Previsto (E) / Vl. Total (H) change.
#    Each entry is: row, new D, new E, new H
$updates = @(
    @(5,   0,                    1306.56,             1306.56),
    @(6,   0,                    955.47,              955.47),
    @(8,   0,                    0,                   0),
    @(15,  0,                    883.67,              883.67),
    @(17,  0,                    1027.29,             1027.29),
    @(33,  0,                    318.49,              318.49),
    @(43,  0,                    443.53,              443.53),
    @(49,  0,                    631.54999999999995,  631.54999999999995),
    @(51,  0,                    887.56,              887.56),
    @(52,  0,                    25492.78,            25492.78),
    @(55,  0,                    250.34,              250.34),
    @(57,  0,                    721.99,              721.99),
    @(58,  0,                    0.54,                0.54),
    @(60,  -9719.7999999999993,  10785.26,            1065.46),
    @(99,  0,                    1531.75,             1531.75),
    @(101, -67.069999999999993,  46.99,               -20.079999999999998),
    @(102, 0,                    1617.68,             1617.68),
    @(103, 0,                    1807.28,             1807.28),
    @(104, 0,                    4068.77,             4068.77),
    @(105, 0,                    565,                 565),
    @(107, 0,                    1022.38,             1022.38),
    @(108, 0,                    2896.19,             2896.19),
    @(109, 0,                    272.36,              272.36),
    @(110, 0,                    3128.51,             3128.51),
    @(112, 0,                    53.14,               53.14),
    @(113, 0,                    71.81,               71.81),
    @(114, 0,                    127.62,              127.62),
    @(115, 0,                    171.69,              171.69),
    @(118, 0,                    916.85,              916.85),
    @(120, 0,                    894.72,              894.72),
    @(130, 0,                    225.1,               225.1),
    @(132, 0,                    978.67,              978.67),
    @(138, 0,                    944.27,              944.27),
    @(143, 0,                    16108.72,            16108.72),
    @(148, 0,                    21214.85,            21214.85),
    @(158, 0,                    82.83,               82.83),
    @(165, 0,                    9.07,                9.07),
    @(173, 0,                    2048.67,             2048.67),
    @(189, 0,                    928.35,              928.35),
    @(205, 0,                    62102.46,            62102.46),
    @(230, 0,                    692.38,              692.38),
    @(231, 0,                    744.79,              744.79),
    @(232, 0,                    37.15,               37.15),
    @(235, 0,                    515.91999999999996,  515.91999999999996),
    @(245, 0,                    3818.68,             3818.68),
    @(249, 0,                    932.01,              932.01),
    @(255, 0,                    73.18,               73.18),
    @(264, 0,                    2743.66,             2743.66),
    @(265, 0,                    1855.09,             1855.09),
    @(270, 0,                    0.95,                0.95),
    @(271, 0,                    1214.81,             1214.81),
    @(273, 0,                    1556.64,             1556.64)
)

foreach ($item in $updates) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value = $item[1]
    $ws.Cells.Item($r, 5).Value = $item[2]
    $ws.Cells.Item($r, 8).Value = $item[3]
}

Write-Output "applied date bump + $($updates.Count) balance updates"
